# Auto-generated script applying the cryptos.xlsx price/volume/coin-order update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.63%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.21%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.184"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.27%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05732"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.91%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.553"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.86%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.104"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.14%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8588"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.80%"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8665"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.60%"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01026"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.53%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1367"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.67%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07093"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.99%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03015"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.80%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09386"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.09%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001536"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.35%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006064"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.45%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.007489"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "5,223.89%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.494"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.49%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.187"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-5.70%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.35%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03322"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.56%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1290"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.24%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.490"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.82%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04140"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.72%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.54%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.001228"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "1.03%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.004990"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "11.47%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "2.63%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03756"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.97%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005816"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.60%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1070"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.08%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002102"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.62%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009460"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.59%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.87%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.08%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05706"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-43.52%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002284"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-9.53%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.08%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.08%"
